$d = $word.ActiveDocument

$d.Content.Find.Execute("422÷9=46, 8", $true, $false, $false, $false, $false, $true, 1, $false, "968÷5=193, 3", 2) | Out-Null
$d.Content.Find.Execute("355÷6=59, 1", $true, $false, $false, $false, $false, $true, 1, $false, "829÷5=165, 4", 2) | Out-Null
$d.Content.Find.Execute("206÷2=103, 0", $true, $false, $false, $false, $false, $true, 1, $false, "295÷8=36, 7", 2) | Out-Null
$d.Content.Find.Execute("403÷8=50, 3", $true, $false, $false, $false, $false, $true, 1, $false, "987÷8=123, 3", 2) | Out-Null
$d.Content.Find.Execute("915÷2=457, 1", $true, $false, $false, $false, $false, $true, 1, $false, "585÷5=117, 0", 2) | Out-Null
$d.Content.Find.Execute("311÷5=62, 1", $true, $false, $false, $false, $false, $true, 1, $false, "533÷7=76, 1", 2) | Out-Null
$d.Content.Find.Execute("695÷2=347, 1", $true, $false, $false, $false, $false, $true, 1, $false, "361÷2=180, 1", 2) | Out-Null
$d.Content.Find.Execute("622÷8=77, 6", $true, $false, $false, $false, $false, $true, 1, $false, "352÷7=50, 2", 2) | Out-Null
$d.Content.Find.Execute("766÷7=109, 3", $true, $false, $false, $false, $false, $true, 1, $false, "901÷6=150, 1", 2) | Out-Null
$d.Content.Find.Execute("612÷7=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "809÷3=269, 2", 2) | Out-Null
$d.Content.Find.Execute("245÷2=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "163÷2=81, 1", 2) | Out-Null
$d.Content.Find.Execute("165÷3=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "251÷6=41, 5", 2) | Out-Null
$d.Content.Find.Execute("140÷8=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "575÷7=82, 1", 2) | Out-Null
$d.Content.Find.Execute("745÷4=186, 1", $true, $false, $false, $false, $false, $true, 1, $false, "269÷4=67, 1", 2) | Out-Null
$d.Content.Find.Execute("896÷4=224, 0", $true, $false, $false, $false, $false, $true, 1, $false, "160÷6=26, 4", 2) | Out-Null
$d.Content.Find.Execute("645÷8=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "279÷2=139, 1", 2) | Out-Null
$d.Content.Find.Execute("203÷4=50, 3", $true, $false, $false, $false, $false, $true, 1, $false, "987÷5=197, 2", 2) | Out-Null
$d.Content.Find.Execute("191÷5=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "550÷3=183, 1", 2) | Out-Null
$d.Content.Find.Execute("896÷6=149, 2", $true, $false, $false, $false, $false, $true, 1, $false, "671÷9=74, 5", 2) | Out-Null
$d.Content.Find.Execute("143÷4=35, 3", $true, $false, $false, $false, $false, $true, 1, $false, "491÷7=70, 1", 2) | Out-Null
$d.Content.Find.Execute("691÷6=115, 1", $true, $false, $false, $false, $false, $true, 1, $false, "521÷3=173, 2", 2) | Out-Null
$d.Content.Find.Execute("165÷4=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "554÷6=92, 2", 2) | Out-Null
$d.Content.Find.Execute("983÷7=140, 3", $true, $false, $false, $false, $false, $true, 1, $false, "736÷7=105, 1", 2) | Out-Null
$d.Content.Find.Execute("738÷3=246, 0", $true, $false, $false, $false, $false, $true, 1, $false, "538÷7=76, 6", 2) | Out-Null
$d.Content.Find.Execute("222÷8=27, 6", $true, $false, $false, $false, $false, $true, 1, $false, "481÷8=60, 1", 2) | Out-Null
